$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.000000000000000005310557155020105
$ws.Range("C3").Value = 0.4223910921055155
$ws.Range("C4").Value = 0.00000000000000001647883696604348
$ws.Range("C5").Value = 0.000000000000000022026897714741
$ws.Range("C6").Value = 0.3266462031418655
$ws.Range("C7").Value = 0.0000000000000000008476572151579813
$ws.Range("C8").Value = 0.000000000000000005424406555728697
$ws.Range("C9").Value = 0.02173475375986822
$ws.Range("C10").Value = 0.00000000000000000196004155244827
$ws.Range("C11").Value = 0.000000000000000006794645933682366
$ws.Range("C12").Value = 0.0000000000000000020495369838691
$ws.Range("C13").Value = 0.2292279509927508
